$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 502.54544  # was 517.53845
$ws.Range("J38").Value = 980  # was 871.4286
$ws.Range("L38").Value = 2940  # was 2614.2858
$ws.Range("N38").Value = -3684  # was -3358.2858
$ws.Range("H42").Value = 176  # was 189.5
$ws.Range("I42").Value = 100  # was 75
$ws.Range("J42").Value = 226.66667  # was 246.75
$ws.Range("K42").Value = 300  # was 225
$ws.Range("L42").Value = 680.00001  # was 740.25
$ws.Range("M42").Value = -70  # was 5
$ws.Range("N42").Value = -1140.00001  # was -1200.25
$ws.Range("H52").Value = 8500  # was 8600
$ws.Range("J52").Value = 8500  # was 8600
$ws.Range("L52").Value = 25500  # was 25800
$ws.Range("N52").Value = -25820  # was -26120
$ws.Range("H58").Value = 2563  # was 1347.4166
$ws.Range("I58").Value = 813.3333  # was 416.9
$ws.Range("J58").Value = 3875.25  # was 6000
$ws.Range("K58").Value = 2439.9999  # was 1250.7
$ws.Range("L58").Value = 11625.75  # was 18000
$ws.Range("M58").Value = -2289.9999  # was -1100.7
$ws.Range("N58").Value = -11925.75  # was -18300
$ws.Range("H61").Value = 523.3333  # was 707.5
$ws.Range("I61").Value = 585  # was 1015
$ws.Range("K61").Value = 1755  # was 3045
$ws.Range("M61").Value = -1583  # was -2873
$ws.Range("H64").Value = 3979.8667  # was 3956.1875
$ws.Range("I64").Value = 3633.111  # was 3633.2222
$ws.Range("J64").Value = 4500  # was 4371.4287
$ws.Range("K64").Value = 3633.111  # was 3633.2222
$ws.Range("L64").Value = 4500  # was 4371.4287
$ws.Range("M64").Value = -3385.111  # was -3385.2222
$ws.Range("N64").Value = -4996  # was -4867.4287
$ws.Range("H67").Value = 3979.8667  # was 3956.1875
$ws.Range("I67").Value = 3633.111  # was 3633.2222
$ws.Range("J67").Value = 4500  # was 4371.4287
$ws.Range("K67").Value = 3633.111  # was 3633.2222
$ws.Range("L67").Value = 4500  # was 4371.4287
$ws.Range("M67").Value = -2775.111  # was -2775.2222
$ws.Range("N67").Value = -6216  # was -6087.4287
$ws.Range("H98").Value = 1119.6  # was 999.7143
$ws.Range("I98").Value = 1119.6  # was 999.7143
$ws.Range("K98").Value = 1119.6  # was 999.7143
$ws.Range("M98").Value = 378.4000000000001  # was 498.2857
$ws.Range("H113").Value = 31253708  # was 30306664
$ws.Range("I113").Value = 71431780  # was 66669748
$ws.Range("K113").Value = 71431780  # was 66669748
$ws.Range("M113").Value = -71428526  # was -66666494
$ws.Range("H122").Value = 1119.6  # was 999.7143
$ws.Range("I122").Value = 1119.6  # was 999.7143
$ws.Range("K122").Value = 3358.8  # was 2999.1429
$ws.Range("M122").Value = -908.7999999999997  # was -549.1428999999998
$ws.Range("H129").Value = 1127.8704  # was 1105.8474
$ws.Range("J129").Value = 1226.4894  # was 1192.0193
$ws.Range("L129").Value = 3679.4682  # was 3576.0579
$ws.Range("N129").Value = -13679.4682  # was -13576.0579

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4647.383  # was 5052.07
$ws.Range("I32").Value = 4807.4546  # was 5258.5
$ws.Range("K32").Value = 4807.4546  # was 5258.5
$ws.Range("M32").Value = -4520.4546  # was -4971.5
$ws.Range("H45").Value = 3892.6667  # was 3603.3157
$ws.Range("I45").Value = 4052.6667  # was 3751.9
$ws.Range("J45").Value = 3732.6667  # was 3438.2222
$ws.Range("K45").Value = 4052.6667  # was 3751.9
$ws.Range("L45").Value = 3732.6667  # was 3438.2222
$ws.Range("M45").Value = -3675.6667  # was -3374.9
$ws.Range("N45").Value = -4486.6667  # was -4192.2222
$ws.Range("H132").Value = 14762.737  # was 15570.056
$ws.Range("I132").Value = 1499.75  # was 1584.3334
$ws.Range("K132").Value = 4499.25  # was 4753.0002
$ws.Range("M132").Value = -1969.25  # was -2223.0002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1268.0344  # was 1477
$ws.Range("I99").Value = 1016.6957  # was 1199.55
$ws.Range("J99").Value = 2231.5  # was 2401.8333
$ws.Range("K99").Value = 1016.6957  # was 1199.55
$ws.Range("L99").Value = 2231.5  # was 2401.8333
$ws.Range("M99").Value = 481.3043  # was 298.45
$ws.Range("N99").Value = -5227.5  # was -5397.8333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3193.0645  # was 3116.7188
$ws.Range("I31").Value = 2161.6316  # was 2091.05
$ws.Range("K31").Value = 2161.6316  # was 2091.05
$ws.Range("M31").Value = -1866.6316  # was -1796.05
$ws.Range("H34").Value = 3193.0645  # was 3116.7188
$ws.Range("I34").Value = 2161.6316  # was 2091.05
$ws.Range("K34").Value = 2161.6316  # was 2091.05
$ws.Range("M34").Value = -1959.6316  # was -1889.05
$ws.Range("H52").Value = 37690  # was 37740
$ws.Range("J52").Value = 37690  # was 37740
$ws.Range("L52").Value = 37690  # was 37740
$ws.Range("N52").Value = -38278  # was -38328
$ws.Range("H111").Value = 0  # was 35500.332
$ws.Range("J111").Value = 0  # was 35500.332
$ws.Range("L111").Value = 0  # was 35500.332
$ws.Range("N111").ClearContents()  # was -43680.332
$ws.Range("H122").Value = 1373.8462  # was 1584.9048
$ws.Range("I122").Value = 1638.4667  # was 1914
$ws.Range("J122").Value = 1013  # was 1146.1111
$ws.Range("K122").Value = 4915.4001  # was 5742
$ws.Range("L122").Value = 3039  # was 3438.3333
$ws.Range("M122").Value = -2465.4001  # was -3292
$ws.Range("N122").Value = -7939  # was -8338.3333
$ws.Range("H134").Value = 780.54285  # was 825.1875
$ws.Range("I134").Value = 693.6667  # was 716.64
$ws.Range("J134").Value = 1073.75  # was 1212.8572
$ws.Range("K134").Value = 2081.0001  # was 2149.92
$ws.Range("L134").Value = 3221.25  # was 3638.5716
$ws.Range("M134").Value = 453.9998999999998  # was 385.0799999999999
$ws.Range("N134").Value = -8291.25  # was -8708.571599999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 135.72728  # was 166.08333
$ws.Range("J18").Value = 0  # was 500
$ws.Range("L18").Value = 0  # was 1500
$ws.Range("N18").ClearContents()  # was -1838
$ws.Range("H86").Value = 100000280  # was 500
$ws.Range("I86").Value = 300  # was 0
$ws.Range("J86").Value = 125000270  # was 500
$ws.Range("K86").Value = 900  # was 0
$ws.Range("L86").Value = 375000810  # was 1500
$ws.Range("M86").Value = 286  # was None
$ws.Range("N86").Value = -375003182  # was -3872
$ws.Range("H89").Value = 100000280  # was 500
$ws.Range("I89").Value = 300  # was 0
$ws.Range("J89").Value = 125000270  # was 500
$ws.Range("K89").Value = 2700  # was 0
$ws.Range("L89").Value = 1125002430  # was 4500
$ws.Range("M89").Value = 3228  # was None
$ws.Range("N89").Value = -1125014286  # was -16356
$ws.Range("H131").Value = 791.8  # was 789.5599999999999
$ws.Range("J131").Value = 801.81915  # was 799.43616
$ws.Range("L131").Value = 2405.45745  # was 2398.30848
$ws.Range("N131").Value = -12485.45745  # was -12478.30848

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23000  # was 24000
$ws.Range("J15").Value = 23000  # was 24000
$ws.Range("L15").Value = 23000  # was 24000
$ws.Range("N15").Value = -23576  # was -24576
$ws.Range("H81").Value = 23000  # was 24000
$ws.Range("J81").Value = 23000  # was 24000
$ws.Range("L81").Value = 23000  # was 24000
$ws.Range("N81").Value = -24996  # was -25996
$ws.Range("H84").Value = 23000  # was 24000
$ws.Range("J84").Value = 23000  # was 24000
$ws.Range("L84").Value = 69000  # was 72000
$ws.Range("N84").Value = -78984  # was -81984
$ws.Range("H122").Value = 121214770  # was 83335470
$ws.Range("I122").Value = 111112300  # was 41667750
$ws.Range("K122").Value = 333336900  # was 125003250
$ws.Range("M122").Value = -333334450  # was -125000800
$ws.Range("H132").Value = 71598.625  # was 47816.582
$ws.Range("I132").Value = 14448  # was 7089
$ws.Range("J132").Value = 128749.25  # was 169999.33
$ws.Range("K132").Value = 43344  # was 21267
$ws.Range("L132").Value = 386247.75  # was 509997.99
$ws.Range("M132").Value = -40814  # was -18737
$ws.Range("N132").Value = -391307.75  # was -515057.99

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6100.647  # was 5995
$ws.Range("J40").Value = 7061.3  # was 6801.091
$ws.Range("L40").Value = 7061.3  # was 6801.091
$ws.Range("N40").Value = -7333.3  # was -7073.091
$ws.Range("H61").Value = 5055.136  # was 5248.2383
$ws.Range("I61").Value = 2835.8572  # was 2977.077
$ws.Range("K61").Value = 2835.8572  # was 2977.077
$ws.Range("M61").Value = -2633.8572  # was -2775.077
$ws.Range("H113").Value = 5055.136  # was 5248.2383
$ws.Range("I113").Value = 2835.8572  # was 2977.077
$ws.Range("K113").Value = 2835.8572  # was 2977.077
$ws.Range("M113").Value = -665.8571999999999  # was -807.0770000000002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 32700.5  # was 35135.332
$ws.Range("J27").Value = 32700.5  # was 35135.332
$ws.Range("L27").Value = 32700.5  # was 35135.332
$ws.Range("N27").Value = -32838.5  # was -35273.332
$ws.Range("H81").Value = 62501010  # was 71429670
$ws.Range("J81").Value = 166666930  # was 250000210
$ws.Range("L81").Value = 333333860  # was 500000420
$ws.Range("N81").Value = -333335982  # was -500002542
$ws.Range("H84").Value = 62501010  # was 71429670
$ws.Range("J84").Value = 166666930  # was 250000210
$ws.Range("L84").Value = 1666669300  # was 2500002100
$ws.Range("N84").Value = -1666679908  # was -2500012708
$ws.Range("H107").Value = 9092109  # was 7576872.5
$ws.Range("J107").Value = 15152015  # was 11364184
$ws.Range("L107").Value = 45456045  # was 34092552
$ws.Range("N107").Value = -45459885  # was -34096392
$ws.Range("H113").Value = 2253948.2  # was 2080590.6
$ws.Range("I113").Value = 3016.5  # was 2628.4285
$ws.Range("K113").Value = 9049.5  # was 7885.2855
$ws.Range("M113").Value = -6879.5  # was -5715.2855
